$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.866.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.830.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.49%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.42%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.15%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.35%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3689"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.80%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07181"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8767"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.69%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07850"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.59"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.861.25"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.334"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.389"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "86.86"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.67%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008726"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.903.01"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.47"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.993"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.63%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.978"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.94"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.32%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.22"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.963"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.63"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.33%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.96%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08817"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.130"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7530"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.462"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.131"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.553"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.087"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.14%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.931"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05127"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.900"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.76%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.72%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.319"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4678"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.34%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.27"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.17%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06098"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.50"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.45"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.57%  "
